$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the scores for "Change Password" (row 30) and "Authorization Checks" (row 32)
$ws.Range("C30").Value = 5
$ws.Range("C32").Value = 5

# Scroll/selection state update (matches the recorded sheetView change in the diff:
# topLeftCell moved from A28 to A22, active selection moved from C51 to C29)
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C29").Select()
